{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Replace the date line and each \"a\u00f7b=c, d\" table-cell answer with its\n// updated value, using exact (unique) text matches found via\n// context.document.body.search().\n\nconst replacements = [\n  [\"2025-11-29 Saturday\", \"2025-11-30 Sunday\"],\n  [\"290\u00f75=58, 0\", \"717\u00f75=143, 2\"],\n  [\"262\u00f76=43, 4\", \"146\u00f74=36, 2\"],\n  [\"912\u00f76=152, 0\", \"945\u00f76=157, 3\"],\n  [\"172\u00f73=57, 1\", \"964\u00f77=137, 5\"],\n  [\"416\u00f77=59, 3\", \"493\u00f79=54, 7\"],\n  [\"753\u00f75=150, 3\", \"779\u00f78=97, 3\"],\n  [\"445\u00f76=74, 1\", \"109\u00f73=36, 1\"],\n  [\"938\u00f73=312, 2\", \"647\u00f76=107, 5\"],\n  [\"561\u00f75=112, 1\", \"634\u00f75=126, 4\"],\n  [\"451\u00f77=64, 3\", \"884\u00f72=442, 0\"],\n  [\"639\u00f73=213, 0\", \"568\u00f77=81, 1\"],\n  [\"420\u00f72=210, 0\", \"212\u00f74=53, 0\"],\n  [\"294\u00f76=49, 0\", \"317\u00f74=79, 1\"],\n  [\"425\u00f74=106, 1\", \"744\u00f76=124, 0\"],\n  [\"571\u00f79=63, 4\", \"562\u00f76=93, 4\"],\n  [\"349\u00f78=43, 5\", \"928\u00f78=116, 0\"],\n  [\"729\u00f76=121, 3\", \"923\u00f72=461, 1\"],\n  [\"890\u00f74=222, 2\", \"192\u00f78=24, 0\"],\n  [\"773\u00f72=386, 1\", \"693\u00f72=346, 1\"],\n  [\"689\u00f76=114, 5\", \"930\u00f75=186, 0\"],\n  [\"411\u00f75=82, 1\", \"781\u00f79=86, 7\"],\n  [\"881\u00f75=176, 1\", \"202\u00f76=33, 4\"],\n  [\"581\u00f72=290, 1\", \"954\u00f78=119, 2\"],\n  [\"538\u00f75=107, 3\", \"458\u00f78=57, 2\"],\n  [\"235\u00f78=29, 3\", \"809\u00f78=101, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Replace the date line and each \"a\u00f7b=c, d\" table-cell answer with its\n# updated value using Find/Replace over the whole document content.\n\n$pairs = @(\n    @(\"2025-11-29 Saturday\", \"2025-11-30 Sunday\"),\n    @(\"290\u00f75=58, 0\", \"717\u00f75=143, 2\"),\n    @(\"262\u00f76=43, 4\", \"146\u00f74=36, 2\"),\n    @(\"912\u00f76=152, 0\", \"945\u00f76=157, 3\"),\n    @(\"172\u00f73=57, 1\", \"964\u00f77=137, 5\"),\n    @(\"416\u00f77=59, 3\", \"493\u00f79=54, 7\"),\n    @(\"753\u00f75=150, 3\", \"779\u00f78=97, 3\"),\n    @(\"445\u00f76=74, 1\", \"109\u00f73=36, 1\"),\n    @(\"938\u00f73=312, 2\", \"647\u00f76=107, 5\"),\n    @(\"561\u00f75=112, 1\", \"634\u00f75=126, 4\"),\n    @(\"451\u00f77=64, 3\", \"884\u00f72=442, 0\"),\n    @(\"639\u00f73=213, 0\", \"568\u00f77=81, 1\"),\n    @(\"420\u00f72=210, 0\", \"212\u00f74=53, 0\"),\n    @(\"294\u00f76=49, 0\", \"317\u00f74=79, 1\"),\n    @(\"425\u00f74=106, 1\", \"744\u00f76=124, 0\"),\n    @(\"571\u00f79=63, 4\", \"562\u00f76=93, 4\"),\n    @(\"349\u00f78=43, 5\", \"928\u00f78=116, 0\"),\n    @(\"729\u00f76=121, 3\", \"923\u00f72=461, 1\"),\n    @(\"890\u00f74=222, 2\", \"192\u00f78=24, 0\"),\n    @(\"773\u00f72=386, 1\", \"693\u00f72=346, 1\"),\n    @(\"689\u00f76=114, 5\", \"930\u00f75=186, 0\"),\n    @(\"411\u00f75=82, 1\", \"781\u00f79=86, 7\"),\n    @(\"881\u00f75=176, 1\", \"202\u00f76=33, 4\"),\n    @(\"581\u00f72=290, 1\", \"954\u00f78=119, 2\"),\n    @(\"538\u00f75=107, 3\", \"458\u00f78=57, 2\"),\n    @(\"235\u00f78=29, 3\", \"809\u00f78=101, 1\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
